$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 389; this shifts rows 389..483 down to 390..484
$ws.Rows("389:389").Insert()

# Populate the newly inserted row 389 with its data.
$ws.Range("A389").Value = 9
$ws.Range("B389").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C389").Value = "Metropolitana"
$ws.Range("D389").Value = 44798
$ws.Range("E389").Value = 13
$ws.Range("F389").Value = "Fruta"
$ws.Range("G389").Value = 100108
$ws.Range("H389").Value = "Tropicales y subtropicales"
$ws.Range("I389").Value = 100108002
$ws.Range("J389").Value = "Mango"
$ws.Range("K389").Value = "Sin especificar"
$ws.Range("L389").Value = "Primera"
$ws.Range("M389").Value = 630
$ws.Range("N389").Value = 8500
$ws.Range("O389").Value = 9000
$ws.Range("P389").Value = 8778
$ws.Range("Q389").Value = "$/bandeja 4 kilos"
$ws.Range("R389").Value = "México"
$ws.Range("S389").Value = 2194
$ws.Range("T389").Value = 4
